$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "password" column (column D) entirely — header + all values.
$ws.Range("D1").EntireColumn.Delete()
